$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 56: MI vs KXI (match 47) -- fill in raw fantasy points
$ws.Range("E56").Value = 80
$ws.Range("H56").Value = 40
$ws.Range("K56").Value = 60
$ws.Range("N56").Value = 20
$ws.Range("Q56").Value = 100
$ws.Range("T56").Value = 20

# Row 57: MI vs RCB (match 48) -- fill in raw fantasy points
$ws.Range("E57").Value = 20
$ws.Range("H57").Value = 100
$ws.Range("K57").Value = 40
$ws.Range("N57").Value = 100
$ws.Range("Q57").Value = 60
$ws.Range("T57").Value = 0

# M56 and S56 were overwritten with plain (hardcoded) values instead of formulas
$ws.Range("M56").Value = -22.5
$ws.Range("S56").Value = -22.5

# G57 and M57 were overwritten with plain (hardcoded) values instead of formulas
$ws.Range("G57").Value = 35
$ws.Range("M57").Value = 35

$excel.CalculateFull()
